# "financing features of agents added"
# Adds four new columns (L:O) to Sheet1 describing agents' financing setup:
#   L = days_between_financing, M = financing_period,
#   N = ordering_period,        O = delivery_period
# plus sample values for the first three agent rows (2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
# Written in this order so the new shared-strings are appended in the same
# sequence as the target workbook (financing_period, ordering_period,
# delivery_period, days_between_financing).
$ws.Range("M1").Value = "financing_period"
$ws.Range("N1").Value = "ordering_period"
$ws.Range("O1").Value = "delivery_period"
$ws.Range("L1").Value = "days_between_financing"

# --- Data rows --------------------------------------------------------------
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 90

$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 80

$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 70

# --- Column widths (best-fit-like sizing for the new columns) --------------
$ws.Columns.Item(12).ColumnWidth = 23.857142857142858  # L ~ 24.625
$ws.Columns.Item(13).ColumnWidth = 13.571428571428571  # M ~ 14.25
$ws.Columns.Item(14).ColumnWidth = 12.857142857142858  # N ~ 13.625
$ws.Columns.Item(15).ColumnWidth = 12.428571428571429  # O ~ 13.125

# --- View state: the edit selected the two new "period" columns ------------
$ws.Columns("M").Select()
